# HELP_DESK_QUERY_HISTORY_DATA.xlsx ("AMS" sheet) - Bug fix to the
# Manage Interviewers history log: correct the 06-16 run time, and
# record the two 06-17 sprint history runs that were missing (one of
# which produced a case-count mismatch, so it stays highlighted).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AMS")

# --- Row 9 (2021-06-16 / live_145_hotfix): correct the recorded run time ---
$ws.Range("B9").Value2 = 44363.7008096875

# --- Row 10 (2021-06-17 / 145_data_hstry): fill in the previously-blank row ---
$ws.Range("A10").NumberFormat = "@"
$ws.Range("A10").Value2 = "2021-06-17"
$ws.Range("A10").NumberFormat = "General"

$ws.Range("B10").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B10").Value2 = 44364.55006693287

$ws.Range("C10").Value2 = "145_data_hstry"
$ws.Range("D10").Value2 = 124
$ws.Range("E10").Value2 = 122
$ws.Range("F10").Value2 = 2
$ws.Range("G10").Value2 = 1.51

# --- Row 11 (2021-06-17 / 145_hstry_data): fill in the previously-blank row ---
$ws.Range("A11").NumberFormat = "@"
$ws.Range("A11").Value2 = "2021-06-17"
$ws.Range("A11").NumberFormat = "General"

$ws.Range("B11").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B11").Value2 = 44364.61309472586

$ws.Range("C11").Value2 = "145_hstry_data"
# D11 already carries the "mismatch" highlight style from the template;
# just populate its value, total cases still 124 but only 123 passed.
$ws.Range("D11").Value2 = 124
$ws.Range("E11").Value2 = 123
$ws.Range("F11").Value2 = 1
$ws.Range("G11").Value2 = 1.47
